$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume table cells to reflect the latest scrape.
# For cells whose new text value looks like a plain number (e.g. "130.10"),
# force the Text number format first so Excel keeps storing the original
# string (with trailing zeros / 2-decimal formatting) instead of coercing it
# into a numeric value.

$ws.Range("D2").Value = "62.810.39"
$ws.Range("E2").Value = "  +1.44%  "
$ws.Range("D3").Value = "3.476.85"
$ws.Range("E3").Value = "  +2.02%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "416.23"
$ws.Range("E5").Value = "  +1.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.10"
$ws.Range("E6").Value = "  +1.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.629"
$ws.Range("E7").Value = "  -0.82%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.728"
$ws.Range("E9").Value = "  -0.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.153"
$ws.Range("E10").Value = "  +7.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.52"
$ws.Range("E11").Value = "  -0.84%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.79"
$ws.Range("E12").Value = "  +5.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000226"
$ws.Range("E13").Value = "  +2.39%  "
$ws.Range("D14").Value = "4.030.80"
$ws.Range("E14").Value = "  +2.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.140"
$ws.Range("E15").Value = "  -0.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.58"
$ws.Range("E16").Value = "  -3.12%  "
$ws.Range("D17").Value = "3.472.82"
$ws.Range("E17").Value = "  +1.88%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.69"
$ws.Range("E18").Value = "  +1.24%  "
$ws.Range("E19").Value = "  -0.82%  "
$ws.Range("D20").Value = "62.776.18"
$ws.Range("E20").Value = "  +1.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "465.71"
$ws.Range("E21").Value = "  +4.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "90.51"
$ws.Range("E22").Value = "  -1.76%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.30"
$ws.Range("E23").Value = "  +3.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.28"
$ws.Range("E24").Value = "  +0.78%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.78"
$ws.Range("E26").Value = "  +0.93%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "33.33"
$ws.Range("E27").Value = "  +0.83%  "
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("E29").Value = "  -2.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "12.14"
$ws.Range("E30").Value = "  +1.38%  "
$ws.Range("E31").Value = "  -0.73%  "
$ws.Range("E32").Value = "  -0.88%  "
$ws.Range("E33").Value = "  -1.40%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "41.00"
$ws.Range("E34").Value = "  -4.27%  "
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "58.10"
$ws.Range("E36").Value = "  +8.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0492"
$ws.Range("E37").Value = "  -2.46%  "
$ws.Range("E38").Value = "  +0.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.06"
$ws.Range("E39").Value = "  +3.56%  "
$ws.Range("B40").Value = "WEMIXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.73"
$ws.Range("E40").Value = "  +6.42%  "
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "148.03"
$ws.Range("E41").Value = "  +3.19%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.135"
$ws.Range("E42").Value = "  -0.33%  "
$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.323"
$ws.Range("E43").Value = "  +0.65%  "
$ws.Range("B44").Value = "LidoDAOToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.34"
$ws.Range("E44").Value = "  -1.18%  "
$ws.Range("E45").Value = "  +2.89%  "
$ws.Range("E46").Value = "  +3.08%  "
$ws.Range("D47").Value = "0.0₃0592"
$ws.Range("E47").Value = "  +39.63%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.39"
$ws.Range("E48").Value = "  +11.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "16.38"
$ws.Range("E49").Value = "  -1.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.39"
$ws.Range("E50").Value = "  +0.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.142"
$ws.Range("E51").Value = "  -2.91%  "
